# Update the attendance figures (column C) for the FC Barcelona fixture list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "C1" = "80"
    "C2" = "72"
    "C3" = "72"
    "C4" = "63"
    "C5" = "52"
    "C6" = "62"
    "C7" = "53"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the literal to be stored as text (matches the source data, which
    # keeps these numeric-looking values as strings) without leaving a
    # lingering custom number-format style behind.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
